$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:50 PM"

# --- Industry Analysis sheet: update "1 Year" (column F) values ---
$ia = $wb.Worksheets.Item("Industry Analysis")
$ia.Cells.Item(2, 6).Value = 21.0016
$ia.Cells.Item(3, 6).Value = -16.2396
$ia.Cells.Item(4, 6).Value = 27.1317
$ia.Cells.Item(5, 6).Value = -50.6494
$ia.Cells.Item(6, 6).Value = 53.2813
$ia.Cells.Item(7, 6).Value = -8.106199999999999
$ia.Cells.Item(8, 6).Value = -9.552099999999999
$ia.Cells.Item(9, 6).Value = 36.3756
$ia.Cells.Item(10, 6).Value = -6.1314
$ia.Cells.Item(11, 6).Value = 31.9081
$ia.Cells.Item(12, 6).Value = -18.4955
$ia.Cells.Item(13, 6).Value = 14.0155
$ia.Cells.Item(14, 6).Value = -36.0718
$ia.Cells.Item(15, 6).Value = -0.1622
$ia.Cells.Item(16, 6).Value = 0.1459
$ia.Cells.Item(17, 6).Value = -22.0012
$ia.Cells.Item(18, 6).Value = 1.0561
$ia.Cells.Item(19, 6).Value = -27.708
$ia.Cells.Item(20, 6).Value = 47.7309
$ia.Cells.Item(21, 6).Value = 12.0959
$ia.Cells.Item(22, 6).Value = 95.1491
$ia.Cells.Item(23, 6).Value = -50.2657
$ia.Cells.Item(24, 6).Value = -13.3427
$ia.Cells.Item(25, 6).Value = -9.9316
$ia.Cells.Item(26, 6).Value = 5.8244
$ia.Cells.Item(27, 6).Value = -32.7692
$ia.Cells.Item(28, 6).Value = -24.8224
$ia.Cells.Item(29, 6).Value = -18.4191
$ia.Cells.Item(30, 6).Value = 25.8569
$ia.Cells.Item(31, 6).Value = 58.4712
$ia.Cells.Item(32, 6).Value = -3.3862
$ia.Cells.Item(33, 6).Value = -6.3282
$ia.Cells.Item(34, 6).Value = 27.7203
$ia.Cells.Item(35, 6).Value = 4.4873
$ia.Cells.Item(36, 6).Value = -4.9458
$ia.Cells.Item(37, 6).Value = 3.6074
$ia.Cells.Item(38, 6).Value = -23.3973
$ia.Cells.Item(39, 6).Value = 8.7355
$ia.Cells.Item(40, 6).Value = -5.8541
$ia.Cells.Item(41, 6).Value = -8.3934
$ia.Cells.Item(42, 6).Value = 20.3818
$ia.Cells.Item(43, 6).Value = 14.3164
$ia.Cells.Item(44, 6).Value = -12.6846
$ia.Cells.Item(45, 6).Value = 28.4075
$ia.Cells.Item(46, 6).Value = -1.1135
$ia.Cells.Item(47, 6).Value = -37.1997
$ia.Cells.Item(48, 6).Value = -29.8569
$ia.Cells.Item(49, 6).Value = -27.5511
$ia.Cells.Item(50, 6).Value = -49.7478
$ia.Cells.Item(51, 6).Value = -51.8002
$ia.Cells.Item(52, 6).Value = -38.5254
$ia.Cells.Item(53, 6).Value = -12.4886
$ia.Cells.Item(54, 6).Value = -5.0725
$ia.Cells.Item(55, 6).Value = -17.7445
$ia.Cells.Item(56, 6).Value = -26.636
$ia.Cells.Item(57, 6).Value = -29.3361
$ia.Cells.Item(58, 6).Value = -11.9574
$ia.Cells.Item(59, 6).Value = -24.5687
$ia.Cells.Item(60, 6).Value = -12.3
$ia.Cells.Item(61, 6).Value = -10.9446
$ia.Cells.Item(62, 6).Value = -17.1229
$ia.Cells.Item(63, 6).Value = -9.5038
$ia.Cells.Item(64, 6).Value = 54.2749
$ia.Cells.Item(65, 6).Value = -43.4736
$ia.Cells.Item(66, 6).Value = 13.2687
$ia.Cells.Item(67, 6).Value = 12.7149
$ia.Cells.Item(68, 6).Value = 24.8057
$ia.Cells.Item(69, 6).Value = -17.0328
$ia.Cells.Item(70, 6).Value = -6.8927
$ia.Cells.Item(71, 6).Value = 13.6034
$ia.Cells.Item(72, 6).Value = 3.9995
$ia.Cells.Item(73, 6).Value = -16.226
$ia.Cells.Item(74, 6).Value = -16.2448
$ia.Cells.Item(75, 6).Value = 28.6924
$ia.Cells.Item(76, 6).Value = 48.9752
